$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prioritise content")

# Update header text (E2, I2) per the lesson-plan template wording update
$ws.Range("E2").Value = "How important is the task to the learning outcome of the topic ? (3,2,1,0)"
$ws.Range("I2").Value = "allocated time (min)"

# Re-enter the H and K column formulas across the table so they share one
# formula definition (as Excel does when you fill a formula down a range)
$ws.Range("H3:H13").Formula = "=SUM(D3:G3)"
$ws.Range("K3:K13").Formula = "=IF(H3<=6,""Good to Know"",""Must to Know"")"

# Update the sheet view: clear the scrolled top-left cell and move the
# active selection to I2
$ws.Range("I2").Select() | Out-Null
